$d = $word.ActiveDocument

# Locate the anchor paragraph: the empty paragraph (w14:paraId="2F8AF7E0")
# that sits 6 paragraphs before the end of the document, right after the
# "Voeding (batterij)" paragraph and before the trailing run of empty
# paragraphs. We find it positionally by counting from the end, which is
# robust to content earlier in the document.
$count = $d.Paragraphs.Count
$anchorIndex = $count - 5
$anchor = $d.Paragraphs.Item($anchorIndex)

# Sanity check: the anchor paragraph must be empty.
if ($anchor.Range.Text.Trim().Length -ne 0) {
    throw "Anchor paragraph $anchorIndex is not empty (text: '$($anchor.Range.Text)')"
}

# Create a new empty paragraph right after the anchor; this new paragraph's
# Range is where we inject the full multi-paragraph OOXML fragment.
$anchor.Range.InsertParagraphAfter()
$target = $d.Paragraphs.Item($anchorIndex + 1).Range

$xmlFragment = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:u w:val="single"/>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:u w:val="single"/>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>De puzzel:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve">De vuilnisbakken zullen geopend worden via een code die via vorige opdrachten wordt gevonden. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>De spelers zullen via een scanner moeten zien welke stukken afval die gevonden of verdiend kunnen worden in de vuilnisbak moeten. Uiteindelijk zal het gewicht van de verschillende vuilnisbakken een code vormen.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>Implementatie:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve">Er zal 1 grote doos hout gemaakt worden die dienst zal doen als alle vuilnisbakken. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve">In deze bak zal alle elektronica van de vuilnisbakken zitten. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve">Aan de bovenkant van deze bak zullen 3 bakjes zijn (PMD, papier en karton, Restafval). Wanneer de spelers een stuk afval willen weg gooien moeten ze dit leggen in het bakje dat ze denken waar het in past. De </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>rfid</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve"> scanner van dit bakje zal dit stuk vuilnis scannen en bepalen of dit al dan niet correct is. Als dit correct is zal een </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>servo</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve"> motor aangestuurd worden en zal het stuk vuilnis in de bak vallen. De </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>gewichtsensor</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve"> kijkt hoeveel dit weegt en zal dit op het scherm laten zien.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>Indien het fout is zal dit aangegeven worden en zal er energie uit de buffer gaan.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve">De scanner zal in het begin van het spel de correcte </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>rfid</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>values</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve"> krijgen van de vuilnisbak. Deze </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>rfid</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve"> scanner zal bestaan uit een lcd, een knop en een </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>rfid</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>-scanner. De lcd zal laten zien of het een correct stuk afval is.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve">Voor de </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>rfid</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve"> scanners gebruiken we de HW-147 (</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>PN532</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t xml:space="preserve"> module)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="nl-NL"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($xmlFragment)

Write-Output "Inserted puzzel/implementatie section after paragraph $anchorIndex"
